$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats/styles) from row 5 down into the new
# row 6 so the new row's date/percent cells pick up the same style indices
# (s="1" for the date column, s="2" for the percentage columns) instead of
# Excel minting brand-new number-format entries.
$ws.Range("A5:W5").Copy() | Out-Null
$ws.Range("A6:W6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the scan results that were missing from row 5 (X5/Y5)
$ws.Range("X5").Value = -1.4100040000000149
$ws.Range("Y5").Value = "Down"

# Append the new scan result row (row 6)
$ws.Range("A6").Value = 42647.885416666664
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "Neutral"
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 14605
$ws.Range("F6").Value = 809
$ws.Range("G6").Value = 54
$ws.Range("H6").Value = 44
$ws.Range("I6").Value = 61
$ws.Range("J6").Value = 38
$ws.Range("K6").Value = 30684
$ws.Range("L6").Value = 161
$ws.Range("M6").Value = 132
$ws.Range("N6").Value = 43
$ws.Range("O6").Value = 27
$ws.Range("P6").Value = "Named"
$ws.Range("Q6").Value = 60.94594728999143
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0.11890000000000001
$ws.Range("T6").Value = 0.0080000000000000002
$ws.Range("U6").Value = 5.99
$ws.Range("V6").Value = "N/A"
$ws.Range("W6").Value = 0
